$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old hyperlinks on A2 and B2 before rewriting cell content.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()

# Write the new search-related header/value columns first, then the
# credential columns, so the shared-string table is built up in the same
# order as the authored workbook.
$ws.Range("C1").Value = "searchString"
$ws.Range("D1").Value = "searchProductString"
$ws.Range("E1").Value = "minPriceValue"
$ws.Range("F1").Value = "maxPriceValue"
$ws.Range("C2").Value = "amazon"
$ws.Range("D2").Value = "dell computers"

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "dummy@email.com"
$ws.Range("B2").Value = "dummyPass"
$ws.Range("B2").Style = "Normal"

$ws.Range("E2").Value = 20000
$ws.Range("F2").Value = 30000

# Re-create the hyperlink for the (now) dummy email address only; the
# password cell no longer carries a hyperlink.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:dummy@email.com")
$ws.Range("A2").HorizontalAlignment = -4131

# New columns get explicit, content-sized widths.
$ws.Columns.Item(3).ColumnWidth = 11.85546875
$ws.Columns.Item(4).ColumnWidth = 19
$ws.Columns.Item(5).ColumnWidth = 14.140625
$ws.Columns.Item(6).ColumnWidth = 14.42578125

# Move the active selection down to A3, like in the edited workbook.
$ws.Range("A3").Select()
